# Auto-generated COM-interop script to refresh Goblin Profits market data
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 685.1177
$ws.Range("I8").Value = 685.1177
$ws.Range("K8").Value = 2055.3531
$ws.Range("M8").Value = -1916.3531
$ws.Range("H11").Value = 117.166664
$ws.Range("I11").Value = 117.166664
$ws.Range("K11").Value = 117.166664
$ws.Range("M11").Value = 22.833336
$ws.Range("H31").Value = 433
$ws.Range("I31").Value = 433
$ws.Range("K31").Value = 1299
$ws.Range("M31").Value = -1069
$ws.Range("H40").Value = 2996.6667
$ws.Range("I40").Value = 993.3333
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 993.3333
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -818.3333
$ws.Range("N40").Value = -5350
$ws.Range("H42").Value = 58824468
$ws.Range("I42").Value = 71429624
$ws.Range("K42").Value = 214288872
$ws.Range("M42").Value = -214288642
$ws.Range("H93").Value = 65537
$ws.Range("J93").Value = 65537
$ws.Range("L93").Value = 65537
$ws.Range("N93").Value = -70529
$ws.Range("H100").Value = 6094.68
$ws.Range("I100").Value = 4040.4443
$ws.Range("K100").Value = 4040.4443
$ws.Range("M100").Value = -3499.4443
$ws.Range("H113").Value = 3923.4614
$ws.Range("I113").Value = 3923.4614
$ws.Range("K113").Value = 3923.4614
$ws.Range("M113").Value = -669.4614000000001
$ws.Range("H116").Value = 4197.6
$ws.Range("J116").Value = 3997
$ws.Range("L116").Value = 3997
$ws.Range("N116").Value = -10881
$ws.Range("H137").Value = 2500
$ws.Range("I137").Value = 2500
$ws.Range("K137").Value = 7500
$ws.Range("M137").Value = -4950
$ws.Range("H138").Value = 5304.9375
$ws.Range("I138").Value = 1075.9286
$ws.Range("K138").Value = 3227.7858
$ws.Range("M138").Value = 1912.2142
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 5000
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").Value = ""
$ws.Range("H61").Value = 5060.4165
$ws.Range("I61").Value = 4594.1904
$ws.Range("K61").Value = 4594.1904
$ws.Range("M61").Value = -4382.1904
$ws.Range("H102").Value = 3384.4644
$ws.Range("I102").Value = 1580.3182
$ws.Range("K102").Value = 1580.3182
$ws.Range("M102").Value = 41.68180000000007
$ws.Range("H122").Value = 13891773
$ws.Range("I122").Value = 22224170
$ws.Range("K122").Value = 66672510
$ws.Range("M122").Value = -66670060
$ws.Range("H136").Value = 5060.4165
$ws.Range("I136").Value = 4594.1904
$ws.Range("K136").Value = 13782.5712
$ws.Range("M136").Value = -11232.5712
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3038.05
$ws.Range("I99").Value = 2515.353
$ws.Range("J99").Value = 6000
$ws.Range("K99").Value = 2515.353
$ws.Range("L99").Value = 6000
$ws.Range("M99").Value = -1017.353
$ws.Range("N99").Value = -8996
$ws.Range("H125").Value = 75000
$ws.Range("J125").Value = 75000
$ws.Range("L125").Value = 75000
$ws.Range("N125").Value = -84840
$ws.Range("H134").Value = 2917.3333
$ws.Range("I134").Value = 3039
$ws.Range("K134").Value = 9117
$ws.Range("M134").Value = -6582
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5496.8945
$ws.Range("I31").Value = 2374
$ws.Range("J31").Value = 10850.429
$ws.Range("K31").Value = 2374
$ws.Range("L31").Value = 10850.429
$ws.Range("M31").Value = -2079
$ws.Range("N31").Value = -11440.429
$ws.Range("H34").Value = 5496.8945
$ws.Range("I34").Value = 2374
$ws.Range("J34").Value = 10850.429
$ws.Range("K34").Value = 2374
$ws.Range("L34").Value = 10850.429
$ws.Range("M34").Value = -2172
$ws.Range("N34").Value = -11254.429
$ws.Range("H68").Value = 59999.5
$ws.Range("J68").Value = 60000
$ws.Range("L68").Value = 60000
$ws.Range("N68").Value = -61498
$ws.Range("H71").Value = 59999.5
$ws.Range("J71").Value = 60000
$ws.Range("L71").Value = 180000
$ws.Range("N71").Value = -187488
$ws.Range("H81").Value = 74665
$ws.Range("J81").Value = 74665
$ws.Range("L81").Value = 74665
$ws.Range("N81").Value = -76661
$ws.Range("H84").Value = 74665
$ws.Range("J84").Value = 74665
$ws.Range("L84").Value = 223995
$ws.Range("N84").Value = -233979
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 225.86667
$ws.Range("I12").Value = 114
$ws.Range("J12").Value = 266.54544
$ws.Range("K12").Value = 342
$ws.Range("L12").Value = 799.63632
$ws.Range("M12").Value = -169
$ws.Range("N12").Value = -1145.63632
$ws.Range("H59").Value = 866.6667
$ws.Range("J59").Value = 2000
$ws.Range("L59").Value = 6000
$ws.Range("N59").Value = -7080
$ws.Range("H92").Value = 3600.3333
$ws.Range("J92").Value = 2520.4
$ws.Range("L92").Value = 7561.200000000001
$ws.Range("N92").Value = -10057.2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 288.48148
$ws.Range("I97").Value = 291.8846
$ws.Range("K97").Value = 291.8846
$ws.Range("M97").Value = 204.1154
$ws.Range("H122").Value = 8312.232
$ws.Range("I122").Value = 8763.235000000001
$ws.Range("J122").Value = 6608.4443
$ws.Range("K122").Value = 26289.705
$ws.Range("L122").Value = 19825.3329
$ws.Range("M122").Value = -23839.705
$ws.Range("N122").Value = -24725.3329
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6726.846
$ws.Range("I40").Value = 3991.6667
$ws.Range("K40").Value = 3991.6667
$ws.Range("M40").Value = -3855.6667
$ws.Range("H46").Value = 4075.7693
$ws.Range("J46").Value = 4165.4165
$ws.Range("L46").Value = 4165.4165
$ws.Range("N46").Value = -4541.4165
$ws.Range("H93").Value = 6133.933
$ws.Range("I93").Value = 2500
$ws.Range("J93").Value = 7455.364
$ws.Range("K93").Value = 2500
$ws.Range("L93").Value = 7455.364
$ws.Range("M93").Value = -1252
$ws.Range("N93").Value = -9951.364
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 7645.4287
$ws.Range("I132").Value = 7908.2104
$ws.Range("J132").Value = 5149
$ws.Range("K132").Value = 23724.6312
$ws.Range("L132").Value = 15447
$ws.Range("M132").Value = -21194.6312
$ws.Range("N132").Value = -20507
